$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Key Stage 4 (KS4) destinations  - provisional"
$ws.Range("A13").Value = "Key Stage 5 (KS5) destinations - provisional"
$ws.Range("D12").Value = "February 2023 - revision"
$ws.Range("D13").Value = "February 2023 - revision"

$ws.Range("D14").Select()
